$wb = $excel.ActiveWorkbook

# This script applies the numeric corrections produced by the scheduled
# market-data refresh run across all profession sheets (ALC, ARM, BSM, CRP,
# CUL, GSM, LTW, WVR). Each Range(...).Value assignment mirrors one changed
# cell from the diff; two cells (WVR!M64, WVR!M67) are cleared because the
# refreshed HQ/NQ prices converged and the profit-breakdown cell no longer
# applies for those rows.

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 55556630
$ws.Range("I11").Value = 55556630
$ws.Range("K11").Value = 55556630
$ws.Range("M11").Value = -55556490
$ws.Range("H39").Value = 8333397.5
$ws.Range("I39").Value = 10000060
$ws.Range("K39").Value = 30000180
$ws.Range("M39").Value = -29999884
$ws.Range("H98").Value = 1025.2115
$ws.Range("I98").Value = 908.05884
$ws.Range("K98").Value = 908.05884
$ws.Range("M98").Value = 589.94116
$ws.Range("H100").Value = 3251.2856
$ws.Range("I100").Value = 2320
$ws.Range("J100").Value = 3949.75
$ws.Range("K100").Value = 2320
$ws.Range("L100").Value = 3949.75
$ws.Range("M100").Value = -1779
$ws.Range("N100").Value = -5031.75
$ws.Range("H101").Value = 397
$ws.Range("J101").Value = 397
$ws.Range("L101").Value = 1191
$ws.Range("N101").Value = -4435
$ws.Range("H122").Value = 1025.2115
$ws.Range("I122").Value = 908.05884
$ws.Range("K122").Value = 2724.17652
$ws.Range("M122").Value = -274.17652
$ws.Range("H132").Value = 3194.3674
$ws.Range("I132").Value = 3314.7908
$ws.Range("K132").Value = 9944.3724
$ws.Range("M132").Value = -7414.3724

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3498.641
$ws.Range("I32").Value = 2683.8
$ws.Range("K32").Value = 2683.8
$ws.Range("M32").Value = -2396.8
$ws.Range("H61").Value = 3822.2727
$ws.Range("I61").Value = 3654.5
$ws.Range("K61").Value = 3654.5
$ws.Range("M61").Value = -3442.5
$ws.Range("H74").Value = 1335.7587
$ws.Range("I74").Value = 1294.963
$ws.Range("K74").Value = 1294.963
$ws.Range("M74").Value = -420.963
$ws.Range("H77").Value = 1335.7587
$ws.Range("I77").Value = 1294.963
$ws.Range("K77").Value = 6474.815
$ws.Range("M77").Value = -2106.815
$ws.Range("H132").Value = 3011.2122
$ws.Range("I132").Value = 2786.862
$ws.Range("K132").Value = 8360.585999999999
$ws.Range("M132").Value = -5830.585999999999
$ws.Range("H136").Value = 3822.2727
$ws.Range("I136").Value = 3654.5
$ws.Range("K136").Value = 10963.5
$ws.Range("M136").Value = -8413.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 40227.8
$ws.Range("I75").Value = 716.6667
$ws.Range("J75").Value = 99494.5
$ws.Range("K75").Value = 716.6667
$ws.Range("L75").Value = 99494.5
$ws.Range("M75").Value = 219.3333
$ws.Range("N75").Value = -101366.5
$ws.Range("H78").Value = 40227.8
$ws.Range("I78").Value = 716.6667
$ws.Range("J78").Value = 99494.5
$ws.Range("K78").Value = 2150.0001
$ws.Range("L78").Value = 298483.5
$ws.Range("M78").Value = 2529.9999
$ws.Range("N78").Value = -307843.5
$ws.Range("H126").Value = 115000
$ws.Range("J126").Value = 115000
$ws.Range("L126").Value = 115000
$ws.Range("N126").Value = -124880

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2094.875
$ws.Range("I16").Value = 1965.5714
$ws.Range("K16").Value = 1965.5714
$ws.Range("M16").Value = -1678.5714
$ws.Range("H58").Value = 3435.45
$ws.Range("I58").Value = 3167.8
$ws.Range("K58").Value = 3167.8
$ws.Range("M58").Value = -2964.8
$ws.Range("H113").Value = 2094.875
$ws.Range("I113").Value = 1965.5714
$ws.Range("K113").Value = 1965.5714
$ws.Range("M113").Value = 204.4286
$ws.Range("H116").Value = 88659
$ws.Range("J116").Value = 88659
$ws.Range("L116").Value = 88659
$ws.Range("N116").Value = -97837
$ws.Range("H120").Value = 35983
$ws.Range("J120").Value = 33975
$ws.Range("L120").Value = 33975
$ws.Range("N120").Value = -41233
$ws.Range("H121").Value = 45000
$ws.Range("J121").Value = 45000
$ws.Range("L121").Value = 45000
$ws.Range("N121").Value = -47620
$ws.Range("H134").Value = 3533.182
$ws.Range("I134").Value = 2858.8
$ws.Range("J134").Value = 4978.2856
$ws.Range("K134").Value = 8576.400000000001
$ws.Range("L134").Value = 14934.8568
$ws.Range("M134").Value = -6041.400000000001
$ws.Range("N134").Value = -20004.8568
$ws.Range("H136").Value = 3435.45
$ws.Range("I136").Value = 3167.8
$ws.Range("K136").Value = 9503.400000000001
$ws.Range("M136").Value = -6953.400000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 10000
$ws.Range("I58").Value = 10000
$ws.Range("K58").Value = 30000
$ws.Range("M58").Value = -29872
$ws.Range("H97").Value = 862
$ws.Range("J97").Value = 1097.4
$ws.Range("L97").Value = 3292.2
$ws.Range("N97").Value = -4284.200000000001
$ws.Range("H103").Value = 327.7143
$ws.Range("I103").Value = 335.66666
$ws.Range("J103").Value = 280
$ws.Range("K103").Value = 1006.99998
$ws.Range("L103").Value = 840
$ws.Range("M103").Value = -127.9999799999999
$ws.Range("N103").Value = -2598
$ws.Range("H132").Value = 911.1667
$ws.Range("I132").Value = 593.4
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 5340.599999999999
$ws.Range("L132").Value = 22500
$ws.Range("M132").Value = -2810.599999999999
$ws.Range("N132").Value = -27560
$ws.Range("H140").Value = 1390
$ws.Range("I140").Value = 1348.4
$ws.Range("K140").Value = 4045.2
$ws.Range("M140").Value = 1134.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H114").Value = 67995.8
$ws.Range("J114").Value = 46744.75
$ws.Range("L114").Value = 46744.75
$ws.Range("N114").Value = -55422.75
$ws.Range("H122").Value = 1823.1177
$ws.Range("I122").Value = 1443.6666
$ws.Range("K122").Value = 4330.9998
$ws.Range("M122").Value = -1880.9998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9238.657999999999
$ws.Range("I7").Value = 8446.76
$ws.Range("J7").Value = 10761.538
$ws.Range("K7").Value = 8446.76
$ws.Range("L7").Value = 10761.538
$ws.Range("M7").Value = -8334.76
$ws.Range("N7").Value = -10985.538
$ws.Range("H126").Value = 9238.657999999999
$ws.Range("I126").Value = 8446.76
$ws.Range("J126").Value = 10761.538
$ws.Range("K126").Value = 25340.28
$ws.Range("L126").Value = 32284.614
$ws.Range("M126").Value = -22870.28
$ws.Range("N126").Value = -37224.614
$ws.Range("H132").Value = 4566
$ws.Range("I132").Value = 4999
$ws.Range("J132").Value = 3700
$ws.Range("K132").Value = 14997
$ws.Range("L132").Value = 11100
$ws.Range("M132").Value = -12467
$ws.Range("N132").Value = -16160
$ws.Range("H136").Value = 3936.762
$ws.Range("I136").Value = 3713.9375
$ws.Range("K136").Value = 11141.8125
$ws.Range("M136").Value = -8591.8125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 80141
$ws.Range("J46").Value = 80141
$ws.Range("L46").Value = 80141
$ws.Range("N46").Value = -80603
$ws.Range("H64").Value = 99971
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 99971
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 99971
$ws.Range("N64").Value = -100467
$ws.Range("H67").Value = 99971
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 99971
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 99971
$ws.Range("N67").Value = -101687
$ws.Range("H81").Value = 2416.6667
$ws.Range("I81").Value = 1719.6
$ws.Range("J81").Value = 3288
$ws.Range("K81").Value = 3439.2
$ws.Range("L81").Value = 6576
$ws.Range("M81").Value = -2378.2
$ws.Range("N81").Value = -8698
$ws.Range("H84").Value = 2416.6667
$ws.Range("I84").Value = 1719.6
$ws.Range("J84").Value = 3288
$ws.Range("K84").Value = 17196
$ws.Range("L84").Value = 32880
$ws.Range("M84").Value = -11892
$ws.Range("N84").Value = -43488
$ws.Range("H132").Value = 3468.9443
$ws.Range("I132").Value = 3208.0232
$ws.Range("K132").Value = 9624.069600000001
$ws.Range("M132").Value = -7094.069600000001
$ws.Range("H134").Value = 80141
$ws.Range("J134").Value = 80141
$ws.Range("L134").Value = 240423
$ws.Range("N134").Value = -245493
$ws.Range("M64").ClearContents()
$ws.Range("M67").ClearContents()
